$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "contact ba.admission@au.dk" sentence from the qualifying entry
# examination answer (originally in C7).
$ws.Range("C7").Value = 'To apply to Aarhus University, you must hold a qualifying entry examination, which must be passed before 5 July in the year of application (unless you are a paying applicant).

If you have more than one completed entry examination, only the GPA from your first exam can be used in quota 1.
A retake exam cannot be used in quota 1, but additional exams may be used to fulfil specific admission requirements.

A qualifying entry examination never becomes outdated.

Common examples include:

Danish qualifying exams (stx, hf, hhx, htx)

Danish qualifying exam for refugees and immigrants (GIF)

Danish/French Baccalaureate (DFB)

European Baccalaureate (EB)

International Baccalaureate (IB)

Option International du Baccalaureate (OIB)

(I)GCSE, AS-, and A-levels

Most European secondary school exams

American High School Diploma + 1 year of university/college studies or 3 AP tests

Most non-European secondary school degrees + 1 year of university/college studies

Nordic, Faroese, or Greenlandic recognized exams

A full list can be found at the Danish Agency for International Education website.

Applicants with international qualifications should also review Aarhus University’s language requirements.
.'

# Remove the "If you plan to take supplementary courses abroad..." sentence
# from the supplementary-courses answer (originally in C10). The shorter
# text re-wraps to a smaller row height.
$ws.Range("C10").Value = 'Yes. If you do not meet the specific admission requirements, you can supplement your exam.

In Denmark, this is done through Upper Secondary School Supplementing (GS) (available only in Danish).

Outside Denmark, Aarhus University generally recognizes:

International General Certificate of Education (IGCE) from accredited institutions such as Pearson or Cambridge

AS-level ≈ Danish B-level

A-level ≈ Danish A-level

Advanced Placement (AP) Tests issued by the College Board

In Denmark, you can also request an individual assessment of your qualifications at a VUC (Adult Education Centre).'
$ws.Rows.Item(10).RowHeight = 289

# Update the active selection / scroll position to match the edited area.
$ws.Range("C11").Select()
